# Weekly update: re-map the price-report rows (2-17) onto a new set of
# dates/varieties/qualities/prices. The underlying data for columns
# D (Fecha), K (Variedad), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de
# comercializacion), R (Origen), S (Precio $/Kg) and T (Kg / unidad)
# is simply reshuffled between rows; columns A, B, C, E, F, G, H, I, J
# stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values for every row so the rearrangement can be
# applied without one write clobbering a value that still needs to be read.
$snapshot = @{}
for ($r = 2; $r -le 17; $r++) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value()
        K = $ws.Cells.Item($r, 11).Value()
        L = $ws.Cells.Item($r, 12).Value()
        M = $ws.Cells.Item($r, 13).Value()
        N = $ws.Cells.Item($r, 14).Value()
        O = $ws.Cells.Item($r, 15).Value()
        P = $ws.Cells.Item($r, 16).Value()
        Q = $ws.Cells.Item($r, 17).Value()
        R = $ws.Cells.Item($r, 18).Value()
        S = $ws.Cells.Item($r, 19).Value()
        T = $ws.Cells.Item($r, 20).Value()
    }
}

# Destination row -> source row (which row's original data now lands here).
$mapping = @{
    2  = 14
    3  = 15
    4  = 7
    5  = 8
    6  = 4
    7  = 11
    8  = 5
    9  = 2
    10 = 3
    11 = 6
    12 = 9
    13 = 10
    14 = 16
    15 = 17
    16 = 12
    17 = 13
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $src = $snapshot[$srcRow]

    $ws.Cells.Item($destRow, 4).Value = $src.D
    $ws.Cells.Item($destRow, 11).Value = $src.K
    $ws.Cells.Item($destRow, 12).Value = $src.L
    $ws.Cells.Item($destRow, 13).Value = $src.M
    $ws.Cells.Item($destRow, 14).Value = $src.N
    $ws.Cells.Item($destRow, 15).Value = $src.O
    $ws.Cells.Item($destRow, 16).Value = $src.P
    $ws.Cells.Item($destRow, 17).Value = $src.Q
    $ws.Cells.Item($destRow, 18).Value = $src.R
    $ws.Cells.Item($destRow, 19).Value = $src.S
    $ws.Cells.Item($destRow, 20).Value = $src.T
}
